$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConfigCarros")

# Append the new data row (row 2) under the existing header row.
$ws.Range("A2").Value = "teste"
$ws.Range("B2").Value = "teste"
$ws.Range("C2").Value = "teste"
$ws.Range("D2").Value = 2020
$ws.Range("E2").Value = "teste"
$ws.Range("F2").Value = 2000
$ws.Range("G2").Value = "teste"
$ws.Range("H2").Value = "Ativo"
$ws.Range("I2").Value = "17/05/2024 - 11:06:28"

# Re-install the AutoFilter so it spans the newly added row too.
# (AutoFilter() on an already-filtered range toggles it off, so remove the
# old one first, then reapply it across the expanded range.)
$ws.Range("A1:I1").AutoFilter()
$ws.Range("A1:I2").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name in sync with the
# resized AutoFilter range.
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "='ConfigCarros'!`$A`$1:`$I`$2"
